$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-8 (A: Processo, B: Originario, C: Status)
$data = @(
    @("5035699-12.2017.8.21.0001", "0161621-51.2017.8.21.0001", "Digitalizado"),
    @("5000273-46.2011.8.21.0001", "0420275-57.2011.8.21.0001", "Digitalizado"),
    @("5000274-31.2011.8.21.0001", "0413893-48.2011.8.21.0001", "Digitalizado"),
    @("5004401-12.2011.8.21.0001", "0376983-22.2011.8.21.0001", "Digitalizado"),
    @("5000256-73.2012.8.21.0001", "0023506-26.2012.8.21.0001", "Digitalizado"),
    @("5000271-42.2012.8.21.0001", "0351565-48.2012.8.21.0001", "Digitalizado"),
    @("5000293-03.2012.8.21.0001", "0342323-65.2012.8.21.0001", "Digitalizado")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Remove old rows 9-16 which are no longer part of the data set
$ws.Range("A9:C16").ClearContents() | Out-Null
